$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 3 (pushing every
# subsequent data row down by one). Insert a fresh row there and
# populate it with the new week's values.
$ws.Rows(3).Insert()

$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 45245
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 100112026
$ws.Cells.Item(3, 7).Value = "Haba"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 400
$ws.Cells.Item(3, 11).Value = 11000
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 13).Value = 11500
$ws.Cells.Item(3, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(3, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(3, 16).Value = 460
$ws.Cells.Item(3, 17).Value = 25
$ws.Cells.Item(3, 18).Value = "Hortaliza"
